$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "70.056.49"
$ws.Range("E2").Value = "  -0.95%  "

# Row 3
$ws.Range("D3").Value = "2.523.25"
$ws.Range("E3").Value = "  -1.23%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.78"
$ws.Range("E5").Value = "  -0.27%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.06"
$ws.Range("E6").Value = "  -2.04%  "

# Row 7
$ws.Range("E7").Value = "  +0.02%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.519"
$ws.Range("E8").Value = "  +1.54%  "

# Row 9
$ws.Range("D9").Value = "2.522.82"
$ws.Range("E9").Value = "  -1.24%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.162"
$ws.Range("E10").Value = "  -2.40%  "

# Row 11
$ws.Range("E11").Value = "  -0.95%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.355"
$ws.Range("E12").Value = "  +2.76%  "

# Row 13
$ws.Range("E13").Value = "  +2.31%  "

# Row 14
$ws.Range("D14").Value = "2.986.96"
$ws.Range("E14").Value = "  -1.19%  "

# Row 15
$ws.Range("D15").Value = "70.022.95"
$ws.Range("E15").Value = "  -0.75%  "

# Row 16
$ws.Range("E16").Value = "  -2.37%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.05"
$ws.Range("E17").Value = "  -0.03%  "

# Row 18
$ws.Range("D18").Value = "2.528.32"
$ws.Range("E18").Value = "  -1.41%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.46"
$ws.Range("E19").Value = "  -1.75%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.81"
$ws.Range("E20").Value = "  +9.74%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.24"
$ws.Range("E21").Value = "  -3.26%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.95"
$ws.Range("E22").Value = "  -0.42%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.02"
$ws.Range("E23").Value = "  +0.60%  "

# Row 24
$ws.Range("E24").Value = "  +0.00%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.51"
$ws.Range("E25").Value = "  +0.64%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.02"
$ws.Range("E26").Value = "  -2.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.94"
$ws.Range("E27").Value = "  -4.06%  "

# Row 28
$ws.Range("D28").Value = "2.646.95"
$ws.Range("E28").Value = "  -1.50%  "

# Row 29
$ws.Range("E29").Value = "  +0.09%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0904"
$ws.Range("E30").Value = "  -2.87%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.96"
$ws.Range("E31").Value = "  +0.78%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "467.27"
$ws.Range("E32").Value = "  -4.18%  "

# Row 33
$ws.Range("E33").Value = "  -2.19%  "

# Row 34
$ws.Range("E34").Value = "  -1.60%  "

# Row 35
$ws.Range("E35").Value = "  +0.05%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.116"
$ws.Range("E36").Value = "  +0.05%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.02"
$ws.Range("E37").Value = "  -0.17%  "

# Row 38
$ws.Range("E38").Value = "  +1.17%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.61"
$ws.Range("E39").Value = "  -0.74%  "

# Row 41
$ws.Range("E41").Value = "  +0.55%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.320"

# Row 43
$ws.Range("E43").Value = "  -4.01%  "

# Row 44
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.47"
$ws.Range("E44").Value = "  +0.04%  "

# Row 45
$ws.Range("B45").Value = "ImmutableX"
$ws.Range("C45").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.15"
$ws.Range("E45").Value = "  -12.84%  "

# Row 46
$ws.Range("E46").Value = "  -6.74%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "142.79"
$ws.Range("E47").Value = "  -1.57%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.530"
$ws.Range("E48").Value = "  -0.37%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.50"
$ws.Range("E49").Value = "  -1.55%  "

# Row 50
$ws.Range("E50").Value = "  -3.16%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0733"
$ws.Range("E51").Value = "  -0.63%  "
